$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    'freezer dresses for women',
    'freezer pants',
    'freezer pants for men',
    'freezer pants for women',
    'freezer suits',
    'freezer suits for women work',
    'freezer wear for women',
    'frisbee equipment',
    'frisbee gear',
    'frisbee hoop',
    'frisbee key',
    'frisbee pro',
    'frisbee shorts',
    'from fatigued to fantastic',
    'full back support brace',
    'full compression',
    'full compression pants',
    'full figured tops',
    'full house apparel',
    'full leg compression',
    'full leg compression machine',
    'full leggings for women',
    'full length',
    'full length back brace',
    'full length compression',
    'full length leg compression sleeve',
    'full length leggings',
    'full length leggings for women',
    'full length workout leggings for women',
    'full lumbar support',
    'full opaque tights',
    'full pant',
    'full reflective jacket',
    'full size basketball game',
    'full stockings for women',
    'fun dresses for women 2017',
    'fun running gear for women',
    'functional training vest',
    'futbol tennis',
    'futbol women',
    'game day sweat shirts women',
    'game gear compression shorts',
    'game gear shorts',
    'game leggings',
    'games for ladies night',
    'garden ideas for women',
    'garden pants for women',
    'garment tape',
    'gasp knee wraps',
    'gasp tank top men',
    'gear bunch leggings',
    'gear for running in the cold',
    'gear ring',
    'gear sport compression',
    'gear support',
    'gear tights',
    'gear up bike stand',
    'gear ver',
    'generator bike',
    'generator bike light',
    'generator stand',
    'generator warmers',
    'genius joint',
    'gift generator',
    'gifts for male runners',
    'gifts for marathon runners male',
    'gifts for triathletes women',
    'girlfriend leggings',
    'girls compression',
    'glue runner',
    'glute exercises',
    'glute training',
    'glute workout',
    'glutes exercise',
    'glutes workout',
    'go lumbar support',
    'go lumbar support extra comfortable adjustable breathable',
    'goal frame soccer',
    'goal tracking',
    'god leggings',
    'gold running tights',
    'good sport clothing',
    'got apperal',
    'grad walk dress',
    'gradient tights',
    'grafitti leggings',
    'gray capri leggings',
    'green basketball tights',
    'green black tights women',
    'green capri leggings',
    'green compression leggings',
    'green compression pants',
    'green compression tights women',
    'green day band patches',
    'green leggingd',
    'green mens compression leggings',
    'green muscle leggings',
    'green panty hose',
    'green puff balls',
    'green rugby shorts'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
